$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '43.949.81'
$c.Style = "Normal"
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.361.30'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '0.689'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +5.76%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '240.80'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +3.05%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '76.32'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +5.29%  '
$ws.Range('E8').Value = '  +0.02%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.624'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +25.21%  '
$ws.Range('E10').Value = '  +5.31%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '57.37'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +0.87%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '33.10'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +21.21%  '
$ws.Range('E13').Value = '  +19.48%  '
$ws.Range('E14').Value = '  +2.01%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '2.716.24'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +0.43%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '16.89'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +3.74%  '
$ws.Range('E17').Value = '  +6.11%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '2.359.99'
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '43.910.56'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +1.66%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.0000104'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +4.13%  '
$ws.Range('E21').Value = '  +4.90%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '77.38'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +3.73%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '256.44'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +2.13%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '11.32'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +12.81%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +2.74%  '
$ws.Range('E27').Value = '  -6.39%  '
$ws.Range('E28').Value = '  +14.80%  '
$ws.Range('E29').Value = '  +1.66%  '
$ws.Range('E30').Value = '  +3.22%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '174.93'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +1.38%  '
$ws.Range('E32').Value = '  -2.71%  '
$ws.Range('E33').Value = '  +5.61%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '5.30'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +5.65%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0751'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +8.59%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '5.33'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +5.49%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '3.82'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +2.03%  '
$ws.Range('E38').Value = '  +0.14%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '6.47'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.80%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.0276'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +8.48%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '18.96'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.200'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +17.77%  '
$ws.Range('E43').Value = '  -0.01%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '8.91'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('E45').Value = '  +5.55%  '
$ws.Range('E46').Value = '  +4.79%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.19'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.50'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +12.66%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '101.67'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +2.36%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '4.50'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +0.32%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '54.31'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +7.21%  '
